$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "NO.OF.INST" column (old C),
# shifting NO.OF.INST / COMMISSION.TYPE:1 / COMMISSION.AMT:1 (and the row 2
# data under them) two columns to the right.
$ws.Columns("C:D").Insert()

# New header cells for the inserted columns (row 2 under them stays blank).
$ws.Range("C1").Value = "CHEQUE.NUMBER"
$ws.Range("D1").Value = "ORDERING.BANK:1"

# Match the "DEBIT.ACCT.NO" column width on the two new columns.
$ws.Columns("C:D").ColumnWidth = 13.83

# Reflect the selection state recorded for the sheet after the edit.
$ws.Range("D1:D1048576").Select() | Out-Null
